# Applies the "AO updated using shoreline access weights based on roads,
# slope, and MPA or no fishing areas" edit:
#   - Adds a new "Ocean Jobs" block to the bottom of the "LE" sheet.
#   - Adds a new "AO" worksheet (Access / Resource / Status / Trend by
#     region), becoming the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "LE" sheet - append the 2013 Ocean Jobs-by-region block below the
#    existing goal-score table, and widen the two new data columns.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LE")

$ws2.Range("A10").Value = 2013

$ws2.Range("A11").Value = "Region"
$ws2.Range("B11").Value = "Ocean Jobs"

$ws2.Range("A12").Value = "Hawaii"
$ws2.Range("B12").Value = 13576

$ws2.Range("A13").Value = "Maui Nui"
$ws2.Range("B13").Value = 25423

$ws2.Range("A14").Value = "Oahu"
$ws2.Range("B14").Value = 59163

$ws2.Range("A15").Value = "Kauai"
$ws2.Range("B15").Value = 5264

$ws2.Columns.Item(3).ColumnWidth = 14.26
$ws2.Columns.Item(4).ColumnWidth = 17.78

$ws2.Range("C11").Select()

# ---------------------------------------------------------------------
# 2) New "AO" sheet - shoreline access weights by region, placed after
#    "LE" as the last (and active) tab.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$ao = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($sheetCount))
$ao.Name = "AO"

# Region names first (A column) so new shared strings land in the same
# order the source workbook uses: ... Ocean Jobs(27), Kauai & Niʻihau(28),
# Access(29), Resource(30).
$ao.Range("A1").Value = "Region"
$ao.Range("A2").Value = "Hawaiʻi"
$ao.Range("A3").Value = "Maui Nui"
$ao.Range("A4").Value = "Oahu"
$ao.Range("A5").Value = "Kauai & Niʻihau"

$ao.Range("B1").Value = "Access"
$ao.Range("C1").Value = "Resource"
$ao.Range("D1").Value = "Status"
$ao.Range("E1").Value = "Trend"

$ao.Range("B2").Value = 0.64
$ao.Range("C2").Value = 0.66
$ao.Range("D2").Value = 0.65
$ao.Range("E2").Value = 0.01

$ao.Range("B3").Value = 0.56000000000000005
$ao.Range("C3").Value = 0.66
$ao.Range("D3").Value = 0.61
$ao.Range("E3").Value = -0.01

$ao.Range("B4").Value = 0.68
$ao.Range("C4").Value = 0.54
$ao.Range("D4").Value = 0.61
$ao.Range("E4").Value = 0.04

$ao.Range("B5").Value = 0.56000000000000005
$ao.Range("C5").Value = 0.72
$ao.Range("D5").Value = 0.64
$ao.Range("E5").Value = 0.02

$ao.Range("J8").Select()
